$d = $word.ActiveDocument

# Locate the original sentence lead-in that needs to be split into several
# runs (with spell-check proofErr markers around "query" and "aql"). Note
# the non-breaking space the source document uses right before the colon.
$old = "Template de test pour les balises de référence à une variable" + [char]0x00A0 + ": "
$rng = $d.Content
$found = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $startPos = $rng.Start

    # Remove the matched text; this collapses $rng to an insertion point.
    $rng.Delete()

    # Build the replacement run sequence (as WordOpenXML) reproducing the
    # diff: the sentence is now split across several runs, with
    # <w:proofErr> spell-check markers wrapping "query" and "aql".
    $newRuns = '<w:r><w:t>Template de test pour les balises d</w:t></w:r>' + `
               '<w:r><w:t xml:space="preserve">e </w:t></w:r>' + `
               '<w:proofErr w:type="spellStart"/>' + `
               '<w:r><w:t>query</w:t></w:r>' + `
               '<w:proofErr w:type="spellEnd"/>' + `
               '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
               '<w:proofErr w:type="spellStart"/>' + `
               '<w:r><w:t>aql</w:t></w:r>' + `
               '<w:proofErr w:type="spellEnd"/>' + `
               ('<w:r><w:t xml:space="preserve"> ' + [char]0x00A0 + ': </w:t></w:r>')

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:body><w:p>' + $newRuns + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $paraCountBefore = $d.Paragraphs.Count
    $rng.InsertXML($xml)

    # Safety net: in some cases InsertXML on a collapsed range inserts the
    # new content as a brand-new paragraph immediately after the (now
    # empty) original one instead of splicing it in place. Detect that by
    # comparing paragraph counts, and if so, merge the two paragraphs back
    # together by deleting the paragraph mark it introduced, restoring the
    # original single-paragraph shape.
    if ($d.Paragraphs.Count -gt $paraCountBefore) {
        $p = $d.Paragraphs.Item(1)
        $brk = $d.Range($p.Range.End - 1, $p.Range.End)
        $brk.Delete()
    }
}
